$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete "Submitted" / "Expect this to beat N" annotations.
# Clearing every cell that references each shared string lets the exporter
# garbage-collect the unused <si> entries (matching the diff's shared-strings
# table shrink from 28 -> 26 unique strings).
$ws.Range("K9").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("M11").ClearContents()

# Row 16 ("lgbm"): CV GAP score is no longer here; the 0.19789999999999999
# figure and its siblings move down into the new rows below.
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()

# Row 17 ("lgbm120") gains CV GAP (C:E), Leaderboard (F) and the GAP formula (G).
$ws.Range("C17").Value = 0.0228
$ws.Range("C17").NumberFormat = "0.0000"
$ws.Range("D17").Value = 0.0204
$ws.Range("D17").NumberFormat = "0.0000"
$ws.Range("E17").Value = 0.0193
$ws.Range("E17").NumberFormat = "0.0000"
$ws.Range("F17").Value = 0.01795
$ws.Range("F17").NumberFormat = "0.00000"
$ws.Range("G17").Formula = "=E17-F17"
$ws.Range("G17").NumberFormat = "0.00000"

# Insert a new row 18 ("lgbm 0.0226"), inheriting row 17's formatting.
$ws.Rows("18:18").Insert()
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "lgbm 0.0226"
$ws.Range("C18").Value = 0.0251
$ws.Range("D18").Value = 0.0217
$ws.Range("E18").Value = 0.0226
$ws.Range("F18:G18").Clear()

# Two trailing blank rows (19, 20), numbered in column A only.
$ws.Rows("19:20").Insert()
$ws.Range("B19:G20").Clear()
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18

# Final cursor position, as saved in the workbook.
$ws.Range("J14").Select() | Out-Null
